$wb = $excel.ActiveWorkbook
$win = $wb.Windows.Item(1)
Write-Host $win.Width
Write-Host $win.Height
$win.Width = 19500
$win.Height = 7890
Write-Host $win.Width
Write-Host $win.Height
